# Rebuild the "Pairing data" table: the former single wide row (row 3, columns
# B..CT holding a 0-96 permutation) is replaced with the same 97 values laid
# out two-per-row in columns B and C, for rows 2 through 56.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$pairData = @{
    2 = @(5)
    3 = @(63)
    4 = @(96)
    5 = @(16,86)
    6 = @(70,42)
    7 = @(11,54)
    8 = @(90,73)
    9 = @(85,50)
    10 = @(59)
    11 = @(45,27)
    12 = @(71,77)
    13 = @(74,92)
    14 = @(24,1)
    15 = @(29)
    16 = @(80,52)
    17 = @(14,41)
    18 = @(39,19)
    19 = @(49)
    20 = @(51)
    21 = @(47,13)
    22 = @(61,43)
    23 = @(65,12)
    24 = @(55,46)
    25 = @(7,44)
    26 = @(32,10)
    27 = @(81,31)
    28 = @(58,8)
    29 = @(87,18)
    30 = @(0,64)
    31 = @(67,66)
    32 = @(3,34)
    33 = @(9,82)
    34 = @(83,26)
    35 = @(15,37)
    36 = @(53,78)
    37 = @(22,21)
    38 = @(36,93)
    39 = @(88,60)
    40 = @(72,20)
    41 = @(25,17)
    42 = @(75,28)
    43 = @(40,35)
    44 = @(91,84)
    45 = @(48,23)
    46 = @(94,38)
    47 = @(79,76)
    48 = @(62)
    49 = @(57,2)
    50 = @(89,33)
    51 = @(69,56)
    52 = @(6)
    53 = @(95)
    54 = @(30)
    55 = @(68)
    56 = @(4)
}

# Clear everything from column B onward first, since row 3 previously spanned
# all the way out to column CT and those trailing cells must disappear.
$ws.Range("B2:CT98").ClearContents()

foreach ($r in $pairData.Keys) {
    $vals = $pairData[$r]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $ws.Cells.Item($r, 2 + $i).Value = $vals[$i]
    }
}
